$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.753.43'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.481.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.25'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.479.75'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.02%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.95%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.22'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.074.44'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.36'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.10%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.43%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.492.76'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.43%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000176'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.97%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.826.39'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.22'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.65'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.39'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.12%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.620.09'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.90%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.59'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.42%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.14%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.488.48'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.80%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.45'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.87%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.34'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.88%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.25%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '159.80'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.17%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.02%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.45'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.80%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.809'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.25%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.75'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.04%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.01%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.17%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.85'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.414.39'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.57%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.32%  '
